# Auto-generated Excel COM-interop script
# Applies per-cell value updates to the Leve profit tracking sheets
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 7534.95
$ws.Range("J51").Value = 2676.8462
$ws.Range("L51").Value = 2676.8462
$ws.Range("N51").Value = -3644.8462
# Row 62
$ws.Range("H62").Value = 2929.8333
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 3055.8
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 3055.8
$ws.Range("N62").Value = -4303.8
$ws.Range("M62").Value = -1676
# Row 65
$ws.Range("H65").Value = 2929.8333
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 3055.8
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 15279
$ws.Range("N65").Value = -21519
$ws.Range("M65").Value = -8380
# Row 76
$ws.Range("H76").Value = 4333
$ws.Range("I76").Value = 4333
$ws.Range("K76").Value = 4333
$ws.Range("M76").Value = -4018
# Row 79
$ws.Range("H79").Value = 4333
$ws.Range("I79").Value = 4333
$ws.Range("K79").Value = 4333
$ws.Range("M79").Value = -3241
# Row 86
$ws.Range("H86").Value = 93108.45
$ws.Range("I86").Value = 93108.45
$ws.Range("K86").Value = 93108.45
$ws.Range("M86").Value = -91985.45
# Row 89
$ws.Range("H89").Value = 93108.45
$ws.Range("I89").Value = 93108.45
$ws.Range("K89").Value = 465542.25
$ws.Range("M89").Value = -459926.25
# Row 92
$ws.Range("H92").Value = 1273.3846
$ws.Range("I92").Value = 1169.25
$ws.Range("J92").Value = 1440
$ws.Range("K92").Value = 1169.25
$ws.Range("L92").Value = 1440
$ws.Range("M92").Value = 78.75
$ws.Range("N92").Value = -3936
# Row 98
$ws.Range("H98").Value = 413.66666
$ws.Range("I98").Value = 407.5
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 407.5
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 1090.5
$ws.Range("N98").Value = -3496
# Row 122
$ws.Range("H122").Value = 413.66666
$ws.Range("I122").Value = 407.5
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1222.5
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 1227.5
$ws.Range("N122").Value = -6400
# Row 137
$ws.Range("H137").Value = 1813.8334
$ws.Range("I137").Value = 1792.0714
$ws.Range("K137").Value = 5376.2142
$ws.Range("M137").Value = -2826.2142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1913.9
$ws.Range("I61").Value = 1390.85
$ws.Range("J61").Value = 2960
$ws.Range("K61").Value = 1390.85
$ws.Range("L61").Value = 2960
$ws.Range("M61").Value = -1178.85
$ws.Range("N61").Value = -3384
# Row 74
$ws.Range("H74").Value = 965.6445
$ws.Range("I74").Value = 977.35297
$ws.Range("J74").Value = 929.4545000000001
$ws.Range("K74").Value = 977.35297
$ws.Range("L74").Value = 929.4545000000001
$ws.Range("M74").Value = -103.35297
$ws.Range("N74").Value = -2677.4545
# Row 77
$ws.Range("H77").Value = 965.6445
$ws.Range("I77").Value = 977.35297
$ws.Range("J77").Value = 929.4545000000001
$ws.Range("K77").Value = 4886.76485
$ws.Range("L77").Value = 4647.2725
$ws.Range("M77").Value = -518.7648500000005
$ws.Range("N77").Value = -13383.2725
# Row 122
$ws.Range("H122").Value = 1670.8096
$ws.Range("I122").Value = 1518.0588
$ws.Range("J122").Value = 2320
$ws.Range("K122").Value = 4554.1764
$ws.Range("L122").Value = 6960
$ws.Range("M122").Value = -2104.1764
$ws.Range("N122").Value = -11860
# Row 136
$ws.Range("H136").Value = 1913.9
$ws.Range("I136").Value = 1390.85
$ws.Range("J136").Value = 2960
$ws.Range("K136").Value = 4172.549999999999
$ws.Range("L136").Value = 8880
$ws.Range("M136").Value = -1622.549999999999
$ws.Range("N136").Value = -13980

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 16993.334
$ws.Range("I26").Value = 16993.334
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 16993.334
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -16701.334
$ws.Range("N26").ClearContents()
# Row 40
$ws.Range("H40").Value = 25601.75
$ws.Range("J40").Value = 27333.334
$ws.Range("L40").Value = 27333.334
$ws.Range("N40").Value = -27863.334
# Row 134
$ws.Range("H134").Value = 2457.5
$ws.Range("I134").Value = 2170.0435
$ws.Range("K134").Value = 6510.130500000001
$ws.Range("M134").Value = -3975.130500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2241.7334
$ws.Range("I58").Value = 2051.2
$ws.Range("J58").Value = 2622.8
$ws.Range("K58").Value = 2051.2
$ws.Range("L58").Value = 2622.8
$ws.Range("M58").Value = -1848.2
$ws.Range("N58").Value = -3028.8
# Row 81
$ws.Range("H81").Value = 38067.5
$ws.Range("J81").Value = 38067.5
$ws.Range("L81").Value = 38067.5
$ws.Range("N81").Value = -40063.5
# Row 82
$ws.Range("H82").Value = 31544.6
$ws.Range("J82").Value = 31544.6
$ws.Range("L82").Value = 31544.6
$ws.Range("N82").Value = -32266.6
# Row 84
$ws.Range("H84").Value = 38067.5
$ws.Range("J84").Value = 38067.5
$ws.Range("L84").Value = 114202.5
$ws.Range("N84").Value = -124186.5
# Row 85
$ws.Range("H85").Value = 31544.6
$ws.Range("J85").Value = 31544.6
$ws.Range("L85").Value = 31544.6
$ws.Range("N85").Value = -34040.6
# Row 132
$ws.Range("H132").Value = 6544.3125
$ws.Range("I132").Value = 9887.429
$ws.Range("J132").Value = 3944.111
$ws.Range("K132").Value = 29662.287
$ws.Range("L132").Value = 11832.333
$ws.Range("M132").Value = -27132.287
$ws.Range("N132").Value = -16892.333
# Row 134
$ws.Range("H134").Value = 1903.1428
$ws.Range("I134").Value = 1803.6666
$ws.Range("K134").Value = 5410.9998
$ws.Range("M134").Value = -2875.9998
# Row 136
$ws.Range("H136").Value = 2241.7334
$ws.Range("I136").Value = 2051.2
$ws.Range("J136").Value = 2622.8
$ws.Range("K136").Value = 6153.599999999999
$ws.Range("L136").Value = 7868.400000000001
$ws.Range("M136").Value = -3603.599999999999
$ws.Range("N136").Value = -12968.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 577
$ws.Range("I92").Value = 496
$ws.Range("J92").Value = 617.5
$ws.Range("K92").Value = 1488
$ws.Range("L92").Value = 1852.5
$ws.Range("M92").Value = -240
$ws.Range("N92").Value = -4348.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 111227520
$ws.Range("I80").Value = 143005780
$ws.Range("J80").Value = 3640
$ws.Range("K80").Value = 143005780
$ws.Range("L80").Value = 3640
$ws.Range("M80").Value = -143004782
$ws.Range("N80").Value = -5636
# Row 83
$ws.Range("H83").Value = 111227520
$ws.Range("I83").Value = 143005780
$ws.Range("J83").Value = 3640
$ws.Range("K83").Value = 715028900
$ws.Range("L83").Value = 18200
$ws.Range("M83").Value = -715023908
$ws.Range("N83").Value = -28184

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 325440.7
$ws.Range("I55").Value = 494822
$ws.Range("J55").Value = 793.1667
$ws.Range("K55").Value = 494822
$ws.Range("L55").Value = 793.1667
$ws.Range("M55").Value = -494649
$ws.Range("N55").Value = -1139.1667
# Row 100
$ws.Range("H100").Value = 5280
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 5280
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 5280
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -6362
# Row 122
$ws.Range("H122").Value = 2928.4285
$ws.Range("I122").Value = 2874.75
$ws.Range("K122").Value = 8624.25
$ws.Range("M122").Value = -6174.25
# Row 130
$ws.Range("H130").Value = 35469.293
$ws.Range("J130").Value = 35469.293
$ws.Range("L130").Value = 35469.293
$ws.Range("N130").Value = -45509.293
# Row 136
$ws.Range("H136").Value = 2550
$ws.Range("I136").Value = 2300
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 6900
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -4350
$ws.Range("N136").Value = -14000.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 26655.666
$ws.Range("I75").Value = 10118
$ws.Range("J75").Value = 34924.5
$ws.Range("K75").Value = 10118
$ws.Range("L75").Value = 34924.5
$ws.Range("N75").Value = -36796.5
$ws.Range("M75").Value = -9182
# Row 78
$ws.Range("H78").Value = 26655.666
$ws.Range("I78").Value = 10118
$ws.Range("J78").Value = 34924.5
$ws.Range("K78").Value = 30354
$ws.Range("L78").Value = 104773.5
$ws.Range("N78").Value = -114133.5
$ws.Range("M78").Value = -25674
# Row 136
$ws.Range("H136").Value = 2056.261
$ws.Range("I136").Value = 752.1053000000001
$ws.Range("J136").Value = 8251
$ws.Range("K136").Value = 2256.3159
$ws.Range("L136").Value = 24753
$ws.Range("M136").Value = 293.6840999999999
$ws.Range("N136").Value = -29853

Write-Host "Applied scheduled market-data updates to all sheets."